$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "Access to Information:" bullet - fix "accdfess" typo split across runs
#    -> merge into a single clean run "... quickly access and update ..."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Allows doctors and nurses to quickly accdfess and update patient information.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Allows doctors and nurses to quickly access and update patient information.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. "Siemens Healthineers:" - split bold run so "Healthineers" is wrapped in
#    proofErr spellStart/spellEnd (as Word does for words it doesn't recognise)
# ---------------------------------------------------------------------------
$rng = $d.Content.Duplicate
$rng.Find.Execute("Siemens Healthineers:") | Out-Null
$para = $rng.Paragraphs(1)
$prng = $para.Range
$xml = @"
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:pPr>
    <w:numPr>
      <w:ilvl w:val='0'/>
      <w:numId w:val='8'/>
    </w:numPr>
  </w:pPr>
  <w:r>
    <w:rPr><w:b/><w:bCs/></w:rPr>
    <w:t xml:space='preserve'>Siemens </w:t>
  </w:r>
  <w:proofErr w:type='spellStart'/>
  <w:r>
    <w:rPr><w:b/><w:bCs/></w:rPr>
    <w:t>Healthineers</w:t>
  </w:r>
  <w:proofErr w:type='spellEnd'/>
  <w:r>
    <w:rPr><w:b/><w:bCs/></w:rPr>
    <w:t>:</w:t>
  </w:r>
  <w:r>
    <w:t xml:space='preserve'> Provides a system that integrates various aspects of hospital management for better efficiency.</w:t>
  </w:r>
</w:p>
"@
$prng.InsertXML($xml) | Out-Null

# ---------------------------------------------------------------------------
# 3. "Kareo:" - wrap "Kareo" in proofErr spellStart/spellEnd
# ---------------------------------------------------------------------------
$rng = $d.Content.Duplicate
$rng.Find.Execute("Kareo:") | Out-Null
$para = $rng.Paragraphs(1)
$prng = $para.Range
$xml = @"
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:pPr>
    <w:numPr>
      <w:ilvl w:val='0'/>
      <w:numId w:val='10'/>
    </w:numPr>
  </w:pPr>
  <w:proofErr w:type='spellStart'/>
  <w:r>
    <w:rPr><w:b/><w:bCs/></w:rPr>
    <w:t>Kareo</w:t>
  </w:r>
  <w:proofErr w:type='spellEnd'/>
  <w:r>
    <w:rPr><w:b/><w:bCs/></w:rPr>
    <w:t>:</w:t>
  </w:r>
  <w:r>
    <w:t xml:space='preserve'> Offers a cloud-based system for managing patient records, appointments, and billing.</w:t>
  </w:r>
</w:p>
"@
$prng.InsertXML($xml) | Out-Null

# ---------------------------------------------------------------------------
# 4. "Nex tech:" -> "Nextech:" (fix the typo'd space) and wrap in proofErr
# ---------------------------------------------------------------------------
$rng = $d.Content.Duplicate
$rng.Find.Execute("Nex tech:") | Out-Null
$para = $rng.Paragraphs(1)
$prng = $para.Range
$xml = @"
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:pPr>
    <w:numPr>
      <w:ilvl w:val='0'/>
      <w:numId w:val='12'/>
    </w:numPr>
  </w:pPr>
  <w:proofErr w:type='spellStart'/>
  <w:r>
    <w:rPr><w:b/><w:bCs/></w:rPr>
    <w:t>Nextech</w:t>
  </w:r>
  <w:proofErr w:type='spellEnd'/>
  <w:r>
    <w:rPr><w:b/><w:bCs/></w:rPr>
    <w:t>:</w:t>
  </w:r>
  <w:r>
    <w:t xml:space='preserve'> Specialized for certain medical specialties, offering tools tailored to those fields.</w:t>
  </w:r>
</w:p>
"@
$prng.InsertXML($xml) | Out-Null

# ---------------------------------------------------------------------------
# 5. "EClinicalWorks:" -> "eClinicalWorks:" (lowercase the leading E, merge runs)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("EClinicalWorks:", $true, $false, $false, $false, $false, $true, 1, $false, `
    "eClinicalWorks:", 2) | Out-Null

# ---------------------------------------------------------------------------
# 6. "OpenMRS:" - wrap "OpenMRS" in proofErr spellStart/spellEnd
#    (keep the lastRenderedPageBreak marker on the first run)
# ---------------------------------------------------------------------------
$rng = $d.Content.Duplicate
$rng.Find.Execute("OpenMRS:") | Out-Null
$para = $rng.Paragraphs(1)
$prng = $para.Range
$xml = @"
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:pPr>
    <w:numPr>
      <w:ilvl w:val='0'/>
      <w:numId w:val='14'/>
    </w:numPr>
  </w:pPr>
  <w:proofErr w:type='spellStart'/>
  <w:r>
    <w:rPr><w:b/><w:bCs/></w:rPr>
    <w:lastRenderedPageBreak/>
    <w:t>OpenMRS</w:t>
  </w:r>
  <w:proofErr w:type='spellEnd'/>
  <w:r>
    <w:rPr><w:b/><w:bCs/></w:rPr>
    <w:t>:</w:t>
  </w:r>
  <w:r>
    <w:t xml:space='preserve'> An open-source platform used for managing medical records, often used in resource-limited settings.</w:t>
  </w:r>
</w:p>
"@
$prng.InsertXML($xml) | Out-Null

# ---------------------------------------------------------------------------
# 7. "OpenEMR:" - wrap "OpenEMR" in proofErr spellStart/spellEnd
# ---------------------------------------------------------------------------
$rng = $d.Content.Duplicate
$rng.Find.Execute("OpenEMR:") | Out-Null
$para = $rng.Paragraphs(1)
$prng = $para.Range
$xml = @"
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:pPr>
    <w:numPr>
      <w:ilvl w:val='0'/>
      <w:numId w:val='14'/>
    </w:numPr>
  </w:pPr>
  <w:proofErr w:type='spellStart'/>
  <w:r>
    <w:rPr><w:b/><w:bCs/></w:rPr>
    <w:t>OpenEMR</w:t>
  </w:r>
  <w:proofErr w:type='spellEnd'/>
  <w:r>
    <w:rPr><w:b/><w:bCs/></w:rPr>
    <w:t>:</w:t>
  </w:r>
  <w:r>
    <w:t xml:space='preserve'> Another open-source system that provides EHR and practice management features.</w:t>
  </w:r>
</w:p>
"@
$prng.InsertXML($xml) | Out-Null

Write-Output "edits applied"
